# GBDS FILES OCTOBER | GBDS UPDATED
# Applies the October-files update to the Credits workbook:
#  - ROUTE 3 "LESS: CREDIT" table (A27:C28) gets new route/SI/amount data
#  - ROUTE 3 header date (D25) becomes a real date value instead of literal text
#  - Print area is moved from the right-hand (ADD: COLLECTION) table to the
#    left-hand (LESS: CREDIT) table
#  - Selection / active cell moves to the newly edited range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ROUTE 3 / LESS: CREDIT table updates (rows 27-28) ---
$ws.Range("A27").Value = "PLAZA STORE"
$ws.Range("B27").Value = 5837
$ws.Range("C27").Value = 385437

$ws.Range("A28").Value = "INDAY STORE"
$ws.Range("B28").Value = 5838
$ws.Range("C28").Value = 60092

# --- ROUTE 3 header date: was literal text "08/13/2025", now a true date value ---
$ws.Range("D25").Value = 45757

# --- Print area: left-hand table (A25:F36) instead of right-hand table (H25:M36) ---
$ws.PageSetup.PrintArea = '$A$25:$F$36'

# --- Selection moves onto the newly edited PLAZA STORE / INDAY STORE amount range ---
$ws.Activate()
$ws.Range("C28:E28").Select()
